# Fixed misnaming issue with sequences
#
# The run filenames in column A were missing an underscore between "seq"
# and the sequence number (e.g. "seq1_1.xlsx" instead of "seq_1_1.xlsx").
# Insert the missing underscore for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $text = $cell.Value2
    if ($text -ne $null) {
        $cell.Value = $text -replace 'seq(\d)', 'seq_$1'
    }
}

# The B column only ever held stray border/fill formatting with no actual
# content. Drop that leftover formatting: B2 and B4:B11 had no other
# purpose, so remove them outright; B1 and B3 keep their row's font
# styling (bold / red, matching C1 / C3) but lose the unwanted
# border/fill.
$ws.Range("B2").Clear()
$ws.Range("B4:B11").Clear()

$ws.Range("B1").ClearFormats()
$ws.Range("B1").Font.Bold = $true

$ws.Range("B3").ClearFormats()
$ws.Range("B3").Font.Color = $ws.Range("C3").Font.Color

# Restore the view state: zoomed to 80% with the cursor left on B10.
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 80
$ws.Range("B10").Select() | Out-Null
